# Chiffres COVID-19 Valais.xlsx -- update to 15.05.2020 data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Title string in A1 (merged A1:L1)
$ws.Range("A1").Value = "Données COVID-19 Valais 15.05.2020"

# 2) Corrections to a few already-existing rows (new hospital-admission figures
#    for 11-14 May revised upward); H column is a live formula (=G+E) so it
#    recalculates on its own once G changes.
$ws.Range("D75").Value2 = 1
$ws.Range("G75").Value2 = 33
$ws.Range("G76").Value2 = 33
$ws.Range("G77").Value2 = 31
$ws.Range("G78").Value2 = 29

# 3) Row 79 used to be the sheet's trailing "today" placeholder row (date only,
#    plus a few manually-typed figures but no cumulative formulas). Preserve
#    that exact look-and-feel by copying it down to the new last row (80) ...
$ws.Range("A79:L79").Copy()
$ws.Range("A80:L80").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ... then restyle row 79 itself as a normal data row (matching the row above).
$ws.Range("A78:L78").Copy()
$ws.Range("A79:L79").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 4) Row 79 becomes real data for 14.05.2020.
$ws.Range("A79").Value2 = 43965
$ws.Range("B79").Formula = "=B78+C79"
$ws.Range("C79").Value2 = 2
$ws.Range("D79").Value2 = 0
$ws.Range("E79").Value2 = 8
$ws.Range("F79").Value2 = 5
$ws.Range("G79").Value2 = 27
$ws.Range("H79").Formula = "=G79+E79"
$ws.Range("I79").Formula = "=I78+J79"
$ws.Range("J79").Formula = "=K79+L79"
$ws.Range("K79").Value2 = 0
$ws.Range("L79").Value2 = 0

# 5) Row 80 is the new trailing placeholder row for 15.05.2020 (B80/C80 stay
#    blank, like B79/C79 used to be).
$ws.Range("A80").Value2 = 43966
$ws.Range("D80").Value2 = 0
$ws.Range("E80").Value2 = 8
$ws.Range("F80").Value2 = 5
$ws.Range("G80").Value2 = 26
$ws.Range("H80").Formula = "=G80+E80"
$ws.Range("I80").Formula = "=I79+J80"
$ws.Range("J80").Formula = "=K80+L80"
$ws.Range("K80").Value2 = 0
$ws.Range("L80").Value2 = 0

# 6) Scroll the view down to the bottom of the table, matching where the
#    author was working.
$ws.Application.ActiveWindow.ScrollRow = 74
$ws.Range("C74").Select()

$wb.Application.CalculateFull()
